# Switch positions between "Ten Hang" (A) and "Part Number" (B) columns,
# and drop the stale "danh" row (and its corresponding summary rows).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: nhap-linhkien
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Remove the obsolete "danh" row (old row 3) entirely.
$ws1.Rows.Item(3).Delete()

# Swap column A (Ten Hang) and column B (Part Number) for every row,
# including the header row (so the header labels swap too).
$lastRow1 = $ws1.UsedRange.Rows.Count
for ($r = 1; $r -le $lastRow1; $r++) {
    $a = $ws1.Cells.Item($r, 1).Value2
    $b = $ws1.Cells.Item($r, 2).Value2
    $ws1.Cells.Item($r, 1).Value = $b
    $ws1.Cells.Item($r, 2).Value = $a
}

# ---------------------------------------------------------------------
# Sheet 2: xuat-linhkien
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$lastRow2 = $ws2.UsedRange.Rows.Count
for ($r = 1; $r -le $lastRow2; $r++) {
    $a = $ws2.Cells.Item($r, 1).Value2
    $b = $ws2.Cells.Item($r, 2).Value2
    $ws2.Cells.Item($r, 1).Value = $b
    $ws2.Cells.Item($r, 2).Value = $a
}

# ---------------------------------------------------------------------
# Sheet 3: ton-linhkien
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

# Remove the rows that corresponded to the blank-named and "danh" entries
# (old rows 5 then 3 - delete bottom-up so row numbers stay valid).
$ws3.Rows.Item(5).Delete()
$ws3.Rows.Item(3).Delete()
